$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Print Titles: repeat rows 1-3 at the top of every printed page ---
# (writes workbook-level defined name "_xlnm.Print_Titles" for this sheet)
$ws.PageSetup.PrintTitleRows = 'Sheet1!$1:$3'

# --- Selection moved from B6 to B5 ---
$ws.Range("B5").Select()

# --- Page margins switched to the metric (cm) defaults: 1.8 / 1.8 / 1.9 / 1.9 / 0.8 / 0.8 cm ---
# PageSetup margins are expressed in points (1 inch = 72 pt); use the exact
# point values so the stored inch fractions round-trip precisely.
$ws.PageSetup.LeftMargin = 51.0236220472441
$ws.PageSetup.RightMargin = 51.0236220472441
$ws.PageSetup.TopMargin = 53.85826771653544
$ws.PageSetup.BottomMargin = 53.85826771653544
$ws.PageSetup.HeaderMargin = 22.677165354330707
$ws.PageSetup.FooterMargin = 22.677165354330707

# --- Fit-to-page printing: scale 46%, capped to 1 page tall ---
# Setting Zoom first records the 46% scale, then switching on "fit to N page(s)
# tall" (FitToPagesTall) flips the sheet into fit-to-page mode, which is what
# turns on sheetPr/pageSetUpPr@fitToPage="1" while keeping scale="46" stored.
$ws.PageSetup.Zoom = 46
$ws.PageSetup.FitToPagesTall = $false
